$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row below the current row 8 (i.e. at row 9). This inherits the
# formatting of row 8 (style "3"), so the original row 8 content can move there intact
# without Excel manufacturing a brand new composite style.
$ws.Rows.Item(9).Insert()

# Move the original row 8 values ("Appliance Installation" / rate) down into row 9.
$ws.Range("A9").Value = $ws.Range("A8").Value2
$ws.Range("B9").Value = $ws.Range("B8").Value2

# Turn row 8 into the new "AC Unit Cleaning" entry with no special formatting.
$ws.Range("A8:B8").ClearFormats()
$ws.Range("A8").Value = "AC Unit Cleaning"
$ws.Range("B8").Value = "`$80 - `$150"

# The hidden _FilterDatabase defined name range needs to shift down by one row too.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$9:`$B`$73"
    }
}

$ws.Range("B10").Select()
